$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: "First Name" -> "FirstName", "Last Name" -> "LastName"
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"

# Update Id values in rows 2 and 3
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2

# Update the active selection to C1
$ws.Range("C1").Select()
